# Add two new columns (I: "I0", J: "IF") to the header row, matching the
# header/style used by the existing columns (e.g. H1 "IP"), and fill in
# the data row values below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered alignment) from
# the existing "IP" header cell (H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
